$d = $word.ActiveDocument

# Correct the wording: swap "1C00 and 500" -> "500 and 1C00" in the
# GPIO base address explanation sentence.
$d.Content.Find.Execute(
    ", where 1C00 and 500 is the GPIO base address",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ", where 500 and 1C00 is the GPIO base address",
    2)
